# MAR 13 MID - Adding final save on matrix
#
# Adds a new "order" column (C) to the user_list sheet:
#   C1 = "order"
#   C2 = "X – O"
#   C3 = "X – O"
#   C4 = "X – O"
# C3/C4 pick up a distinct (but visually identical) cell style, matching
# the source edit, by round-tripping through the "Normal" cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("user_list")

$ws.Range("C1").Value = "order"

$ws.Range("C2").Value = "X – O"

$ws.Range("C3").Value = "X – O"
$ws.Range("C3").Style = "Normal"

$ws.Range("C4").Value = "X – O"
$ws.Range("C4").Style = "Normal"

# Leave the cursor where the author's final save left it.
[void]$ws.Range("E21").Select()
